# Insert a new row for the "Unknown" county right after "Union" (row 89),
# shifting VanBuren and all subsequent rows down by one, then populate the
# new row's data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(89).Insert()

$ws.Range("A89").Value = "Unknown"
$ws.Range("B89").Value = "0"
$ws.Range("C89").Value = "0"
$ws.Range("D89").Value = "1"
